$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the "ram_en" column (column B) ---
$ws.Range("B1").EntireColumn.Delete()

# --- Remove the "HOLD" row (row 5) ---
$ws.Range("A5").EntireRow.Delete()

# After the above deletions the layout is:
# A=State, B=ram_we, C=head_en, D=tail_en, E=ctr_en, F=addr_src, G=inc_dec

# --- Insert a new column for "outr_en" right after ctr_en (before addr_src) ---
$ws.Range("F1").EntireColumn.Insert()

# New layout:
# A=State, B=ram_we, C=head_en, D=tail_en, E=ctr_en, F=outr_en, G=addr_src, H=inc_dec

# --- Header row ---
$ws.Range("A1").Value = "State"
$ws.Range("B1").Value = "ram_we"
$ws.Range("C1").Value = "head_en"
$ws.Range("D1").Value = "tail_en"
$ws.Range("E1").Value = "ctr_en"
$ws.Range("F1").Value = "outr_en"
$ws.Range("G1").Value = "addr_src"
$ws.Range("H1").Value = "inc_dec"

# --- Row 2: WAIT ---
$ws.Range("A2").Value = "WAIT"
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = "x"
$ws.Range("H2").Value = "x"

# --- Row 3: ENQ ---
$ws.Range("A3").Value = "ENQ"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 0

# --- Row 4: DEQ ---
$ws.Range("A4").Value = "DEQ"
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 1

# --- Fix up view state: selection over the full used range, no stray active cell ---
$ws.Range("A1:H4").Select()
